$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 42608.901574074072
$ws.Range("B5").Value = -30
$ws.Range("C5").Value = 40
$ws.Range("D5").Value = 57
$ws.Range("E5").Value = 26
$ws.Range("F5").Value = 73
$ws.Range("G5").Value = 14776
$ws.Range("H5").Value = 28794
$ws.Range("I5").Value = 2956
$ws.Range("J5").Value = 369
$ws.Range("K5").Value = 524
$ws.Range("L5").Value = 5
$ws.Range("M5").Value = 14
$ws.Range("N5").Value = "Bag"
